$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.05970899999999999
$ws.Range("H2").Value = 0.179127
$ws.Range("I2").Value = 0.003688274646134975
$ws.Range("J2").Value = 0.003688274646134975
$ws.Range("M2").Value = 8.521337333333333
$ws.Range("N2").Value = 25.564012
$ws.Range("O2").Value = 0.2943426187002489
$ws.Range("P2").Value = 0.2943426187002489
$ws.Range("Q2").Value = 0.508800530836
$ws.Range("R2").Value = 4.579204777523999
$ws.Range("S2").Value = 0.001085616417829102
$ws.Range("T2").Value = 0.001085616417829102
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.05970899999999999
$ws.Range("H3").Value = 0.179127
$ws.Range("I3").Value = 0.003688274646134975
$ws.Range("J3").Value = 0.003688274646134975
$ws.Range("O3").Value = 0.1683364841626613
$ws.Range("P3").Value = 0.1683364841626613
$ws.Range("Q3").Value = 0.2909863779809999
$ws.Range("R3").Value = 2.618877401829
$ws.Range("S3").Value = 0.0006208711865566453
$ws.Range("T3").Value = 0.0006208711865566455
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.05970899999999999
$ws.Range("H4").Value = 0.179127
$ws.Range("I4").Value = 0.003688274646134975
$ws.Range("J4").Value = 0.003688274646134975
$ws.Range("O4").Value = 0.5373208971370899
$ws.Range("P4").Value = 0.53732089713709
$ws.Range("Q4").Value = 0.9288126840069998
$ws.Range("R4").Value = 8.359314156062998
$ws.Range("S4").Value = 0.001981787041749227
$ws.Range("T4").Value = 0.001981787041749228
$ws.Range("I5").Value = 0.8850509663933519
$ws.Range("J5").Value = 0.885050966393352
$ws.Range("M5").Value = 8.521337333333333
$ws.Range("N5").Value = 25.564012
$ws.Range("O5").Value = 0.2943426187002489
$ws.Range("P5").Value = 0.2943426187002489
$ws.Range("Q5").Value = 122.0935111190124
$ws.Range("R5").Value = 1098.841600071112
$ws.Range("S5").Value = 0.2605082191314051
$ws.Range("T5").Value = 0.2605082191314052
$ws.Range("I6").Value = 0.8850509663933519
$ws.Range("J6").Value = 0.885050966393352
$ws.Range("O6").Value = 0.1683364841626613
$ws.Range("P6").Value = 0.1683364841626613
$ws.Range("R6").Value = 628.4347554712019
$ws.Range("S6").Value = 0.1489863679874225
$ws.Range("T6").Value = 0.1489863679874226
$ws.Range("I7").Value = 0.8850509663933519
$ws.Range("J7").Value = 0.885050966393352
$ws.Range("O7").Value = 0.5373208971370899
$ws.Range("P7").Value = 0.53732089713709
$ws.Range("S7").Value = 0.4755563792745242
$ws.Range("T7").Value = 0.4755563792745244
$ws.Range("I8").Value = 0.111260758960513
$ws.Range("J8").Value = 0.111260758960513
$ws.Range("M8").Value = 8.521337333333333
$ws.Range("N8").Value = 25.564012
$ws.Range("O8").Value = 0.2943426187002489
$ws.Range("P8").Value = 0.2943426187002489
$ws.Range("Q8").Value = 15.34851350607734
$ws.Range("R8").Value = 138.136621554696
$ws.Range("S8").Value = 0.03274878315101458
$ws.Range("T8").Value = 0.03274878315101459
$ws.Range("I9").Value = 0.111260758960513
$ws.Range("J9").Value = 0.111260758960513
$ws.Range("O9").Value = 0.1683364841626613
$ws.Range("P9").Value = 0.1683364841626613
$ws.Range("S9").Value = 0.01872924498868207
$ws.Range("T9").Value = 0.01872924498868208
$ws.Range("I10").Value = 0.111260758960513
$ws.Range("J10").Value = 0.111260758960513
$ws.Range("O10").Value = 0.5373208971370899
$ws.Range("P10").Value = 0.53732089713709
$ws.Range("S10").Value = 0.05978273082081637
$ws.Range("T10").Value = 0.05978273082081639
